$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database values for rows 2 and 3
$updates = @{
    "D"  = -0.162
    "G"  = -0.239766081871345
    "H"  = -0.239766081871345
    "I"  = -0.5345029239766083
    "J"  = -0.5345029239766083
    "K"  = -0.769
    "L"  = -0.8994152046783626
    "U"  = 0.008999999999999999
    "V"  = 0.004999999999999999
    "W"  = 0.1615546218487395
    "X"  = 0.170317483842711
    "Y"  = -0.008762861993971527
    "Z"  = -0.3280890253261705
    "AA" = 0.1753645433614736
    "AB" = 0.1100380145377869
    "AC" = 0.06532652882368666
    "AD" = 2.04
    "AF" = 2.04
    "AG" = 2.031
    "AH" = 0.53125
    "AI" = -0.6415094339622642
    "AJ" = 0.5301487862176977
    "AK" = -0.6368767638758233
    "AL" = 0.295
    "AM" = 0.295
    "AN" = -10.90909090909091
    "AO" = -1.549152542372882
    "AP" = -10.86096256684492
    "AQ" = -1.549152542372882
}

foreach ($row in 2..3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
